$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 27: cohort 2021, period 4 -> num_customers 43 -> 44, retention_rate recalculated (44/2252)
$ws.Range("C27").Value = 44
$ws.Range("E27").Value = 0.01953818827708704

# Row 36: cohort 2023, period 1 -> num_customers 125 -> 127, retention_rate recalculated (127/1930)
$ws.Range("C36").Value = 127
$ws.Range("E36").Value = 0.06580310880829016

# Row 37: cohort 2023, period 0 -> num_customers 788 -> 794, cohort_size 788 -> 794 (retention_rate stays 1)
$ws.Range("C37").Value = 794
$ws.Range("D37").Value = 794
